# Adapt column header formatting to respective input file names (#7)
# - Rename "_old"/"_new" header suffixes to "_FV2310"/"_FV2404"
# - Freeze the header row
# - Wrap the data range in an Excel Table ("Table1")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseHeaders = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

# Columns A-J (1-10): "<name>_old" -> "<name>_FV2310"
for ($i = 0; $i -lt $baseHeaders.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = "$($baseHeaders[$i])_FV2310"
}

# Column K (11) "diff" stays untouched

# Columns L-U (12-21): "<name>_new" -> "<name>_FV2404"
for ($i = 0; $i -lt $baseHeaders.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = "$($baseHeaders[$i])_FV2404"
}

# Freeze the header row (row 1)
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# Convert the data range into a native Excel Table
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U57"), $null, 1)
$tbl.Name = "Table1"

Write-Host "header renaming, freeze pane and table creation complete"
